# Adapt column header formatting to respective input file names (#7)
#
# - Rename the "_old" / "_new" header-name suffixes to "_FV2404" / "_FV2410"
#   (the AHB format-version identifiers the headers now reference).
# - Turn the sheet's used range into an Excel Table ("Table1").
# - Freeze the header row (pane split after row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "_old"/"_new" header suffixes to "_FV2404"/"_FV2410" ---

$oldHeaders = @{
    "A1" = "Segmentname_FV2404"
    "B1" = "Segmentgruppe_FV2404"
    "C1" = "Segment_FV2404"
    "D1" = "Datenelement_FV2404"
    "E1" = "Segment ID_FV2404"
    "F1" = "Code_FV2404"
    "G1" = "Qualifier_FV2404"
    "H1" = "Beschreibung_FV2404"
    "I1" = "Bedingungsausdruck_FV2404"
    "J1" = "Bedingung_FV2404"
    "L1" = "Segmentname_FV2410"
    "M1" = "Segmentgruppe_FV2410"
    "N1" = "Segment_FV2410"
    "O1" = "Datenelement_FV2410"
    "P1" = "Segment ID_FV2410"
    "Q1" = "Code_FV2410"
    "R1" = "Qualifier_FV2410"
    "S1" = "Beschreibung_FV2410"
    "T1" = "Bedingungsausdruck_FV2410"
    "U1" = "Bedingung_FV2410"
}

foreach ($addr in $oldHeaders.Keys) {
    $ws.Range($addr).Value = $oldHeaders[$addr]
}

# --- 2. Convert the used range A1:U71 into a native Excel table ---

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U71"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ---

$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
